# Swap the contents of columns B and C (header text in row 1, and the
# numeric values in rows 2-13), leaving column A and all styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2

    $bCell.Value2 = $cVal
    $cCell.Value2 = $bVal
}
